$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
# 2016-09-01 04:18:31 -> 2016-09-01 04:19:21 (rows 2 and 4 both shared this text)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 04:19:21"
$wsOverview.Range("G4").Value = "2016-09-01 04:19:21"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
# Priority column (E): ht -> mt (rows 2 and 4)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H): 2016-09-01 04:18:27 -> 2016-09-01 04:19:16 (rows 2 and 4)
$wsZhCn.Range("H2").Value = "2016-09-01 04:19:16"
$wsZhCn.Range("H4").Value = "2016-09-01 04:19:16"
# Correspond Handback DateTime column (K): 2016-09-01 04:18:43 -> 2016-09-01 04:19:34 (rows 2 and 4)
$wsZhCn.Range("K2").Value = "2016-09-01 04:19:34"
$wsZhCn.Range("K4").Value = "2016-09-01 04:19:34"

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
# Priority column (E): ht -> mt (rows 2 and 4) -- same shared string as zh-cn's "ht"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H): 2016-09-01 04:18:31 -> 2016-09-01 04:19:21 (rows 2 and 4)
$wsDeDe.Range("H2").Value = "2016-09-01 04:19:21"
$wsDeDe.Range("H4").Value = "2016-09-01 04:19:21"
# Correspond Handback DateTime column (K): 2016-09-01 04:18:49 -> 2016-09-01 04:19:41 (rows 2 and 4)
$wsDeDe.Range("K2").Value = "2016-09-01 04:19:41"
$wsDeDe.Range("K4").Value = "2016-09-01 04:19:41"
